# The authored change (per the OOXML diff) renames the single worksheet
# from its default Spanish Excel name "Hoja1" to "Data".
# (All other hunks in the diff - fileVersion/rupBuild build stamp,
#  x15ac:absPath (author machine path), xr:revisionPtr document/coauth
#  ids, bookViews window geometry, theme display name, and the
#  baseColWidth/defaultColWidth default-serialization swap - are
#  artifacts of the file being opened/saved by a different Excel
#  build/locale, not content edited through the object model.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Data"
